# Scheduled market-data refresh for the Jenova_Profits workbook.
# Recomputes Leve profitability columns (H-N) per sheet from the
# latest Universalis price snapshot; static values only (no formulas).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 185.71428
$ws.Range("I2").Value2 = 185.71428
$ws.Range("J2").Value2 = 0
$ws.Range("K2").Value2 = 185.71428
$ws.Range("L2").Value2 = 0
$ws.Range("M2").Value2 = -72.71428
$ws.Range("N2").ClearContents()
$ws.Range("H64").Value2 = 5831.92
$ws.Range("I64").Value2 = 4907.5386
$ws.Range("K64").Value2 = 4907.5386
$ws.Range("M64").Value2 = -4659.5386
$ws.Range("H67").Value2 = 5831.92
$ws.Range("I67").Value2 = 4907.5386
$ws.Range("K67").Value2 = 4907.5386
$ws.Range("M67").Value2 = -4049.5386
$ws.Range("I69").Value2 = 13333.333
$ws.Range("J69").Value2 = 17500
$ws.Range("K69").Value2 = 39999.999
$ws.Range("L69").Value2 = 52500
$ws.Range("M69").Value2 = -39125.999
$ws.Range("N69").Value2 = -54248
$ws.Range("I72").Value2 = 13333.333
$ws.Range("J72").Value2 = 17500
$ws.Range("K72").Value2 = 119999.997
$ws.Range("L72").Value2 = 157500
$ws.Range("M72").Value2 = -115631.997
$ws.Range("N72").Value2 = -166236
$ws.Range("H92").Value2 = 735.2105
$ws.Range("I92").Value2 = 335
$ws.Range("J92").Value2 = 1421.2858
$ws.Range("K92").Value2 = 335
$ws.Range("L92").Value2 = 1421.2858
$ws.Range("M92").Value2 = 913
$ws.Range("N92").Value2 = -3917.2858
$ws.Range("H99").Value2 = 1061.875
$ws.Range("I99").Value2 = 489.4
$ws.Range("J99").Value2 = 2016
$ws.Range("K99").Value2 = 1468.2
$ws.Range("L99").Value2 = 6048
$ws.Range("M99").Value2 = 29.80000000000018
$ws.Range("N99").Value2 = -9044
$ws.Range("H127").Value2 = 1646.9166
$ws.Range("I127").Value2 = 1362.7778
$ws.Range("J127").Value2 = 2499.3333
$ws.Range("K127").Value2 = 4088.3334
$ws.Range("L127").Value2 = 7497.999899999999
$ws.Range("M127").Value2 = 871.6665999999996
$ws.Range("N127").Value2 = -17417.9999
$ws.Range("H138").Value2 = 4180.0654
$ws.Range("I138").Value2 = 1184.8334
$ws.Range("J138").Value2 = 6122.919
$ws.Range("K138").Value2 = 3554.5002
$ws.Range("L138").Value2 = 18368.757
$ws.Range("M138").Value2 = 1585.4998
$ws.Range("N138").Value2 = -28648.757
$ws.Range("H141").Value2 = 3720.0625
$ws.Range("I141").Value2 = 3736.6155
$ws.Range("J141").Value2 = 3648.3333
$ws.Range("K141").Value2 = 11209.8465
$ws.Range("L141").Value2 = 10944.9999
$ws.Range("M141").Value2 = -6029.8465
$ws.Range("N141").Value2 = -21304.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 4065.9092
$ws.Range("I61").Value2 = 1346.6
$ws.Range("K61").Value2 = 1346.6
$ws.Range("M61").Value2 = -1134.6
$ws.Range("H97").Value2 = 3459.3333
$ws.Range("I97").Value2 = 3459.3333
$ws.Range("K97").Value2 = 3459.3333
$ws.Range("M97").Value2 = -2963.3333
$ws.Range("H110").Value2 = 1002332.2
$ws.Range("I110").Value2 = 1202398.6
$ws.Range("J110").Value2 = 2000
$ws.Range("K110").Value2 = 1202398.6
$ws.Range("L110").Value2 = 2000
$ws.Range("M110").Value2 = -1200353.6
$ws.Range("N110").Value2 = -6090
$ws.Range("H132").Value2 = 7417.3335
$ws.Range("I132").Value2 = 3389.7693
$ws.Range("J132").Value2 = 12177.182
$ws.Range("K132").Value2 = 10169.3079
$ws.Range("L132").Value2 = 36531.546
$ws.Range("M132").Value2 = -7639.3079
$ws.Range("N132").Value2 = -41591.546
$ws.Range("H136").Value2 = 4065.9092
$ws.Range("I136").Value2 = 1346.6
$ws.Range("K136").Value2 = 4039.8
$ws.Range("M136").Value2 = -1489.8
$ws.Range("H139").Value2 = 50000
$ws.Range("J139").Value2 = 50000
$ws.Range("L139").Value2 = 50000
$ws.Range("N139").Value2 = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value2 = 4487.273
$ws.Range("I134").Value2 = 3856.8572
$ws.Range("J134").Value2 = 5590.5
$ws.Range("K134").Value2 = 11570.5716
$ws.Range("L134").Value2 = 16771.5
$ws.Range("M134").Value2 = -9035.571599999999
$ws.Range("N134").Value2 = -21841.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 741.5714
$ws.Range("I22").Value2 = 815.1667
$ws.Range("J22").Value2 = 300
$ws.Range("K22").Value2 = 815.1667
$ws.Range("L22").Value2 = 300
$ws.Range("M22").Value2 = -465.1667
$ws.Range("N22").Value2 = -1000
$ws.Range("H52").Value2 = 60534.332
$ws.Range("I52").Value2 = 46770.5
$ws.Range("K52").Value2 = 46770.5
$ws.Range("M52").Value2 = -46476.5
$ws.Range("H58").Value2 = 404607.6
$ws.Range("I58").Value2 = 1002096
$ws.Range("J58").Value2 = 6282
$ws.Range("K58").Value2 = 1002096
$ws.Range("L58").Value2 = 6282
$ws.Range("M58").Value2 = -1001893
$ws.Range("N58").Value2 = -6688
$ws.Range("H132").Value2 = 3187.6758
$ws.Range("I132").Value2 = 2743.5833
$ws.Range("J132").Value2 = 4007.5386
$ws.Range("K132").Value2 = 8230.749899999999
$ws.Range("L132").Value2 = 12022.6158
$ws.Range("M132").Value2 = -5700.749899999999
$ws.Range("N132").Value2 = -17082.6158
$ws.Range("H136").Value2 = 404607.6
$ws.Range("I136").Value2 = 1002096
$ws.Range("J136").Value2 = 6282
$ws.Range("K136").Value2 = 3006288
$ws.Range("L136").Value2 = 18846
$ws.Range("M136").Value2 = -3003738
$ws.Range("N136").Value2 = -23946

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value2 = 120.5
$ws.Range("I23").Value2 = 163
$ws.Range("J23").Value2 = 92.166664
$ws.Range("K23").Value2 = 489
$ws.Range("L23").Value2 = 276.499992
$ws.Range("M23").Value2 = -254
$ws.Range("N23").Value2 = -746.499992
$ws.Range("H51").Value2 = 2416.6667
$ws.Range("I51").Value2 = 2250
$ws.Range("K51").Value2 = 6750
$ws.Range("M51").Value2 = -6290
$ws.Range("H109").Value2 = 167562.33
$ws.Range("I109").Value2 = 1074.8
$ws.Range("K109").Value2 = 3224.4
$ws.Range("M109").Value2 = -2184.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2860073.8
$ws.Range("J80").Value2 = 3336670.8
$ws.Range("L80").Value2 = 3336670.8
$ws.Range("N80").Value2 = -3338666.8
$ws.Range("H83").Value2 = 2860073.8
$ws.Range("J83").Value2 = 3336670.8
$ws.Range("L83").Value2 = 16683354
$ws.Range("N83").Value2 = -16693338

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value2 = 100000
$ws.Range("J6").Value2 = 100000
$ws.Range("L6").Value2 = 100000
$ws.Range("N6").Value2 = -100224
$ws.Range("H55").Value2 = 1320.3636
$ws.Range("I55").Value2 = 939.4
$ws.Range("J55").Value2 = 1637.8334
$ws.Range("K55").Value2 = 939.4
$ws.Range("L55").Value2 = 1637.8334
$ws.Range("M55").Value2 = -766.4
$ws.Range("N55").Value2 = -1983.8334
$ws.Range("H61").Value2 = 6576.1113
$ws.Range("I61").Value2 = 5439
$ws.Range("K61").Value2 = 5439
$ws.Range("M61").Value2 = -5237
$ws.Range("H113").Value2 = 6576.1113
$ws.Range("I113").Value2 = 5439
$ws.Range("K113").Value2 = 5439
$ws.Range("M113").Value2 = -3269
$ws.Range("H139").Value2 = 47800
$ws.Range("J139").Value2 = 47800
$ws.Range("L139").Value2 = 47800
$ws.Range("N139").Value2 = -58080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value2 = 47624584
$ws.Range("I122").Value2 = 90913450
$ws.Range("J122").Value2 = 6829.6
$ws.Range("K122").Value2 = 272740350
$ws.Range("L122").Value2 = 20488.8
$ws.Range("M122").Value2 = -272737900
$ws.Range("N122").Value2 = -25388.8
$ws.Range("H141").Value2 = 49999.25
$ws.Range("J141").Value2 = 49999.25
$ws.Range("L141").Value2 = 49999.25
$ws.Range("N141").Value2 = -60359.25
